# Refresh coin Price / Volume(1h) figures (and fix the RenderToken /
# WEMIXToken row ordering) to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.807.12"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "2.081.69"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.26"
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.59"
$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +0.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0787"
$ws.Range("E10").Value = "  -0.43%  "

$ws.Range("E11").Value = "  +3.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.90"
$ws.Range("E12").Value = "  +1.05%  "

$ws.Range("D13").Value = "2.388.22"
$ws.Range("E13").Value = "  -0.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.15"
$ws.Range("E14").Value = "  -0.37%  "

$ws.Range("E15").Value = "  +1.82%  "

$ws.Range("E16").Value = "  +1.23%  "

$ws.Range("D17").Value = "2.077.16"
$ws.Range("E17").Value = "  -0.32%  "

$ws.Range("D18").Value = "37.751.53"
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.15"
$ws.Range("E19").Value = "  -1.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.50"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").Value = "0.0₃0841"
$ws.Range("E21").Value = "  +1.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.70"
$ws.Range("E22").Value = "  +0.42%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("E24").Value = "  -0.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  +1.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.77"
$ws.Range("E26").Value = "  +8.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.05"
$ws.Range("E27").Value = "  +0.97%  "

$ws.Range("E28").Value = "  -1.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.49"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("E30").Value = "  -1.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.122"
$ws.Range("E31").Value = "  +1.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.74"
$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0634"
$ws.Range("E33").Value = "  +0.52%  "

$ws.Range("E34").Value = "  -0.84%  "

$ws.Range("E35").Value = "  -1.77%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.42"
$ws.Range("E36").Value = "  -1.24%  "

$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.82"
$ws.Range("E37").Value = "  -0.78%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.07%  "

$ws.Range("E39").Value = "  +0.27%  "

$ws.Range("E40").Value = "  +9.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.15"
$ws.Range("E41").Value = "  +2.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0972"
$ws.Range("E42").Value = "  -1.19%  "

$ws.Range("E43").Value = "  -0.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.98"
$ws.Range("E44").Value = "  +5.36%  "

$ws.Range("D45").Value = "1.447.20"
$ws.Range("E45").Value = "  -0.83%  "

$ws.Range("E46").Value = "  -1.31%  "

$ws.Range("E47").Value = "  -0.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.11"
$ws.Range("E48").Value = "  -4.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.34"
$ws.Range("E49").Value = "  -1.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.99"
$ws.Range("E50").Value = "  -1.69%  "

$ws.Range("D51").Value = "2.274.27"
$ws.Range("E51").Value = "  -0.09%  "
